$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.842.24"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.305.00"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.80"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.42"
$ws.Range("E6").Value = "  -5.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.503"
$ws.Range("E7").Value = "  -4.03%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.502"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.76"
$ws.Range("E10").Value = "  -3.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  -3.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.661.20"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.49"
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.283.87"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.793"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.758.61"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.66"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.03"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.68"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.96"
$ws.Range("E24").Value = "  -3.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.50"
$ws.Range("E25").Value = "  -3.56%  "
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.68"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.32"
$ws.Range("E29").Value = "  -5.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "164.70"
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.10"
$ws.Range("E31").Value = "  -4.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("E33").Value = "  -4.12%  "
$ws.Range("E34").Value = "  -4.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.49"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.67"
$ws.Range("E36").Value = "  -8.79%  "
$ws.Range("E37").Value = "  -4.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.88"
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.79"
$ws.Range("E39").Value = "  -3.69%  "
$ws.Range("E40").Value = "  -4.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").Value = "  -4.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.60"
$ws.Range("E42").Value = "  +5.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.966.21"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0279"
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.27"
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.19"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.86"
$ws.Range("E47").Value = "  -6.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.21"
$ws.Range("E48").Value = "  -4.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.530.03"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.82"
$ws.Range("E50").Value = "  -4.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.67"
$ws.Range("E51").Value = "  -0.14%  "
